# Update "想去人数" (interest count) and "最低票价" (min ticket price) values
# across the 杭州-漫展信息 workbook sheets, matching the upstream data refresh
# captured in commit 456a3b4 ("Update gh-pages to output generated at 456a3b4").
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 4264
$ws.Range("F6").Value = 415
$ws.Range("F7").Value = 3391
$ws.Range("F8").Value = 951
$ws.Range("F11").Value = 271
$ws.Range("F12").Value = 2294
$ws.Range("F16").Value = 491
$ws.Range("F18").Value = 47
$ws.Range("F19").Value = 9450
$ws.Range("G19").Value = 19.9
$ws.Range("F20").Value = 5883
$ws.Range("F21").Value = 375
$ws.Range("F24").Value = 97
$ws.Range("F25").Value = 821
$ws.Range("F29").Value = 443
$ws.Range("F30").Value = 91
$ws.Range("F33").Value = 4764
$ws.Range("F35").Value = 986
$ws.Range("F36").Value = 119
$ws.Range("F37").Value = 446

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 26

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 8610
$ws.Range("F3").Value = 408
$ws.Range("F4").Value = 1478

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 8610
$ws.Range("F4").Value = 408
$ws.Range("F5").Value = 1478
$ws.Range("F7").Value = 4264
$ws.Range("F9").Value = 415
$ws.Range("F10").Value = 3391
$ws.Range("F11").Value = 951
$ws.Range("F14").Value = 271
$ws.Range("F15").Value = 2294
$ws.Range("F18").Value = 26
$ws.Range("F24").Value = 491
$ws.Range("F26").Value = 47
$ws.Range("F27").Value = 9450
$ws.Range("G27").Value = 19.9
$ws.Range("F30").Value = 375
$ws.Range("F33").Value = 97
$ws.Range("F34").Value = 821
$ws.Range("F37").Value = 443
$ws.Range("F38").Value = 91
$ws.Range("F42").Value = 4764
$ws.Range("F44").Value = 986
$ws.Range("F45").Value = 446

